# Populate the four new event rows (315-318) that were added to the
# "Tabelle1" sheet, each with a Date / Event / Location / Stadt / Link.
# The Link column carries a real hyperlink whose visible text is styled
# like the sheet's other Instagram links (underlined, blue, Calibri 11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-EventRow($Row, $DateSerial, $EventText, $LocationText, $CityText, $LinkUrl) {

    $aRef = "A$Row"
    $bRef = "B$Row"
    $cRef = "C$Row"
    $dRef = "D$Row"
    $eRef = "E$Row"

    # Date (keeps the existing dd.mm.yy date-number style already on column A)
    $ws.Range($aRef).Value = $DateSerial

    # Plain text columns - switch their number format to Text ("@") so they
    # match the rest of the table's style instead of staying "General".
    $ws.Range($bRef).Value = $EventText
    $ws.Range($bRef).NumberFormat = "@"

    $ws.Range($cRef).Value = $LocationText
    $ws.Range($cRef).NumberFormat = "@"

    $ws.Range($dRef).Value = $CityText
    $ws.Range($dRef).NumberFormat = "@"

    # Link column: write the URL as the cell text, switch to Text format,
    # then colour/underline every character of it (mirrors how the other
    # rows keep a plain cell style with the link formatting carried on the
    # shared-string run rather than on the cell itself).
    $ws.Range($eRef).Value = $LinkUrl
    $ws.Range($eRef).NumberFormat = "@"
    $linkChars = $ws.Range($eRef).Characters(1, $LinkUrl.Length)
    $linkChars.Font.Underline = $true
    $linkChars.Font.ColorIndex = 4

    # Register the actual hyperlink relationship so the cell is clickable.
    $ws.Hyperlinks.Add($ws.Range($eRef), $LinkUrl, "", "", $LinkUrl) | Out-Null

    # Re-apply the text number format / rich-text run colouring, since
    # attaching the hyperlink can reset the cell-level format.
    $ws.Range($eRef).NumberFormat = "@"
    $linkChars = $ws.Range($eRef).Characters(1, $LinkUrl.Length)
    $linkChars.Font.Underline = $true
    $linkChars.Font.ColorIndex = 4
}

Set-EventRow 315 45760 "AFTERHOUR (6-13Uhr)" "Projekt X" "Bochum" "https://www.instagram.com/reel/DIEPDb9MfVE/?igsh=MWdzNm0zYmZ4MXZrMA=="

Set-EventRow 316 45772 "WYLDHEARTS" "Schrotty" "Köln" "https://www.instagram.com/wyldhearts_?igsh=MWFmY25rN3cyY25mMA=="

Set-EventRow 317 45794 "TECHNOHEADZ" "Elektroküche" "Köln" "https://www.instagram.com/reel/DHG2YbXsRLQ/?igsh=MThka2E5ZmhqbHVoNg=="

Set-EventRow 318 45787 "MATTERMIND" "Essigfabrik & Elektroküche" "Köln" "https://www.instagram.com/reel/DIJd5hyqcWu/?igsh=bTNqYjgwcWlwZ2px"

Write-Output "Added rows 315-318"
